# Updated cryptos list on Sat Nov 18 18:27:49 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "36.626.00"
$ws.Range("E2").Value = "  +0.52%  "
$ws.Range("D3").Value = "1.960.57"
$ws.Range("E3").Value = "  +0.99%  "
$ws.Range("E4").Value = "  -0.16%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "244.33"
$ws.Range("E5").Value = "  +0.46%  "
$ws.Range("E6").Value = "  +0.90%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "58.89"
$ws.Range("E7").Value = "  +2.11%  "
$ws.Range("E8").Value = "  -0.12%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.373"
$ws.Range("E9").Value = "  +1.65%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0804"
$ws.Range("E10").Value = "  -3.32%  "
$ws.Range("E11").Value = "  -0.40%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "22.19"
$ws.Range("E12").Value = "  +3.36%  "
$ws.Range("D13").Value = "2.251.81"
$ws.Range("E13").Value = "  +1.15%  "
$ws.Range("E14").Value = "  -0.55%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "13.69"
$ws.Range("E15").Value = "  +0.67%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.26"
$ws.Range("E16").Value = "  +0.28%  "
$ws.Range("D17").Value = "1.966.22"
$ws.Range("E17").Value = "  +0.27%  "
$ws.Range("D18").Value = "36.538.92"
$ws.Range("E18").Value = "  +0.40%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "69.64"
$ws.Range("E19").Value = "  -0.16%  "
$ws.Range("D20").Value = "0.0₃0858"
$ws.Range("E20").Value = "  -1.13%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "228.59"
$ws.Range("E21").Value = "  -0.32%  "
$ws.Range("E22").Value = "  +1.46%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.999"
$ws.Range("E23").Value = "  +0.02%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.41"
$ws.Range("E24").Value = "  -1.55%  "
$ws.Range("E25").Value = "  +1.62%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.30"
$ws.Range("E26").Value = "  +0.12%  "
$ws.Range("E27").Value = "  +12.49%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "160.21"
$ws.Range("E28").Value = "  -1.23%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "19.34"
$ws.Range("E29").Value = "  -0.20%  "
$ws.Range("E30").Value = "  +1.14%  "
$ws.Range("E31").Value = "  -2.77%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.69"
$ws.Range("E32").Value = "  +0.17%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0616"
$ws.Range("E33").Value = "  -2.27%  "
$ws.Range("E34").Value = "  -0.47%  "
$ws.Range("E35").Value = "  -0.10%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.10"
$ws.Range("E36").Value = "  -0.40%  "
$ws.Range("E37").Value = "  +4.34%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.38"
$ws.Range("E38").Value = "  +13.48%  "
$ws.Range("E39").Value = "  -0.43%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.101"
$ws.Range("E40").Value = "  +4.22%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.90"
$ws.Range("E41").Value = "  -1.47%  "
$ws.Range("E42").Value = "  +1.39%  "
$ws.Range("E43").Value = "  -1.19%  "
$ws.Range("E44").Value = "  -0.13%  "
$ws.Range("D45").Value = "1.359.54"
$ws.Range("E45").Value = "  +0.86%  "
$ws.Range("E46").Value = "  +0.40%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "87.37"
$ws.Range("E47").Value = "  -0.36%  "
$ws.Range("E48").Value = "  -1.13%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.83"
$ws.Range("E49").Value = "  +0.61%  "
$ws.Range("D50").Value = "2.142.06"
$ws.Range("E50").Value = "  +1.06%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "43.39"
$ws.Range("E51").Value = "  -4.92%  "
